$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets updated values and shortened text
$ws.Cells.Item(2, 1).Value = "parsimony. "
$ws.Cells.Item(2, 2).Value = 654
$ws.Cells.Item(2, 3).Value = 669
$ws.Cells.Item(2, 4).Value = 98
$ws.Cells.Item(2, 5).Value = 23

# New row 3: "colossal "
$ws.Cells.Item(3, 1).Value = "colossal "
$ws.Cells.Item(3, 2).Value = 168
$ws.Cells.Item(3, 3).Value = 719
$ws.Cells.Item(3, 4).Value = 75
$ws.Cells.Item(3, 5).Value = 23

# New row 4: "profit "
$ws.Cells.Item(4, 1).Value = "profit "
$ws.Cells.Item(4, 2).Value = 243
$ws.Cells.Item(4, 3).Value = 719
$ws.Cells.Item(4, 4).Value = 48
$ws.Cells.Item(4, 5).Value = 23

# New row 5: "masquerades "
$ws.Cells.Item(5, 1).Value = "masquerades "
$ws.Cells.Item(5, 2).Value = 775
$ws.Cells.Item(5, 3).Value = 785.4
$ws.Cells.Item(5, 4).Value = 122
$ws.Cells.Item(5, 5).Value = 23

# New row 6: "benevolent "
$ws.Cells.Item(6, 1).Value = "benevolent "
$ws.Cells.Item(6, 2).Value = 587
$ws.Cells.Item(6, 3).Value = 968.1999999999999
$ws.Cells.Item(6, 4).Value = 100
$ws.Cells.Item(6, 5).Value = 23

# New row 7: "bequeathed "
$ws.Cells.Item(7, 1).Value = "bequeathed "
$ws.Cells.Item(7, 2).Value = 292
$ws.Cells.Item(7, 3).Value = 993.1999999999999
$ws.Cells.Item(7, 4).Value = 107
$ws.Cells.Item(7, 5).Value = 23
